# chore: update Sheets via scheduled runner
# Refresh market-price derived columns (currentAveragePrice*, LevePrice*, LeveProfit*)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR Leve tables.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
$ws.Range("H34").Value = 23331.666
$ws.Range("I34").Value = 23331.666
$ws.Range("K34").Value = 23331.666
$ws.Range("M34").Value = -23128.666
$ws.Range("H36").Value = 23331.666
$ws.Range("I36").Value = 23331.666
$ws.Range("K36").Value = 23331.666
$ws.Range("M36").Value = -22616.666
$ws.Range("H137").Value = 2793.7646
$ws.Range("I137").Value = 1365
$ws.Range("K137").Value = 4095
$ws.Range("M137").Value = -1545
$ws.Range("H138").Value = 4861.905
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 4861.905
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 14585.715
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -24865.715

$ws = $wb.Worksheets("ARM")
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("H36").Value = 5375
$ws.Range("I36").Value = 2166.6667
$ws.Range("K36").Value = 2166.6667
$ws.Range("M36").Value = -1820.6667
$ws.Range("H125").Value = 155000
$ws.Range("J125").Value = 155000
$ws.Range("L125").Value = 155000
$ws.Range("N125").Value = -164840
$ws.Range("H132").Value = 1849.9846
$ws.Range("I132").Value = 1495.8853
$ws.Range("K132").Value = 4487.6559
$ws.Range("M132").Value = -1957.6559

$ws = $wb.Worksheets("BSM")
$ws.Range("H74").Value = 58000
$ws.Range("J74").Value = 58000
$ws.Range("L74").Value = 58000
$ws.Range("N74").Value = -59872
$ws.Range("H77").Value = 58000
$ws.Range("J77").Value = 58000
$ws.Range("L77").Value = 174000
$ws.Range("N77").Value = -183360
$ws.Range("H139").Value = 99000
$ws.Range("J139").Value = 99000
$ws.Range("L139").Value = 99000
$ws.Range("N139").Value = -109280

$ws = $wb.Worksheets("CRP")
$ws.Range("H22").Value = 687.95654
$ws.Range("I22").Value = 696.5
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 696.5
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -346.5
$ws.Range("N22").Value = -1200
$ws.Range("H31").Value = 2450.7805
$ws.Range("I31").Value = 1512.5555
$ws.Range("J31").Value = 3185.0435
$ws.Range("K31").Value = 1512.5555
$ws.Range("L31").Value = 3185.0435
$ws.Range("M31").Value = -1217.5555
$ws.Range("N31").Value = -3775.0435
$ws.Range("H34").Value = 2450.7805
$ws.Range("I34").Value = 1512.5555
$ws.Range("J34").Value = 3185.0435
$ws.Range("K34").Value = 1512.5555
$ws.Range("L34").Value = 3185.0435
$ws.Range("M34").Value = -1310.5555
$ws.Range("N34").Value = -3589.0435
$ws.Range("H41").Value = 31632.375
$ws.Range("I41").Value = 13264.75
$ws.Range("J41").Value = 50000
$ws.Range("K41").Value = 13264.75
$ws.Range("L41").Value = 50000
$ws.Range("M41").Value = -12836.75
$ws.Range("N41").Value = -50856
$ws.Range("H58").Value = 3042.6072
$ws.Range("I58").Value = 2233
$ws.Range("J58").Value = 3976.7693
$ws.Range("K58").Value = 2233
$ws.Range("L58").Value = 3976.7693
$ws.Range("M58").Value = -2030
$ws.Range("N58").Value = -4382.7693
$ws.Range("H121").Value = 40000
$ws.Range("J121").Value = 40000
$ws.Range("L121").Value = 40000
$ws.Range("N121").Value = -42620
$ws.Range("H132").Value = 2683.7917
$ws.Range("I132").Value = 1731.6471
$ws.Range("J132").Value = 4996.143
$ws.Range("K132").Value = 5194.9413
$ws.Range("L132").Value = 14988.429
$ws.Range("M132").Value = -2664.9413
$ws.Range("N132").Value = -20048.429
$ws.Range("H136").Value = 3042.6072
$ws.Range("I136").Value = 2233
$ws.Range("J136").Value = 3976.7693
$ws.Range("K136").Value = 6699
$ws.Range("L136").Value = 11930.3079
$ws.Range("M136").Value = -4149
$ws.Range("N136").Value = -17030.3079

$ws = $wb.Worksheets("CUL")
$ws.Range("H140").Value = 2688.0605
$ws.Range("I140").Value = 1106.0714
$ws.Range("J140").Value = 3853.7368
$ws.Range("K140").Value = 3318.2142
$ws.Range("L140").Value = 11561.2104
$ws.Range("M140").Value = 1861.7858
$ws.Range("N140").Value = -21921.2104

$ws = $wb.Worksheets("GSM")
$ws.Range("H21").Value = 1685000
$ws.Range("J21").Value = 20000
$ws.Range("L21").Value = 20000
$ws.Range("N21").Value = -20346
$ws.Range("H30").Value = 1685000
$ws.Range("J30").Value = 20000
$ws.Range("L30").Value = 20000
$ws.Range("N30").Value = -20210
$ws.Range("H122").Value = 5324.5293
$ws.Range("J122").Value = 2331.3333
$ws.Range("L122").Value = 6993.999899999999
$ws.Range("N122").Value = -11893.9999
$ws.Range("H132").Value = 14592.45
$ws.Range("I132").Value = 9553.468999999999
$ws.Range("K132").Value = 28660.407
$ws.Range("M132").Value = -26130.407

$ws = $wb.Worksheets("LTW")
$ws.Range("H16").Value = 13199.2
$ws.Range("I16").Value = 27625.25
$ws.Range("J16").Value = 3581.8333
$ws.Range("K16").Value = 27625.25
$ws.Range("L16").Value = 3581.8333
$ws.Range("M16").Value = -27455.25
$ws.Range("N16").Value = -3921.8333
$ws.Range("H22").Value = 1184.3334
$ws.Range("I22").Value = 1221.2
$ws.Range("K22").Value = 1221.2
$ws.Range("M22").Value = -926.2
$ws.Range("H27").Value = 1184.3334
$ws.Range("I27").Value = 1221.2
$ws.Range("K27").Value = 1221.2
$ws.Range("M27").Value = -1114.2
$ws.Range("H34").Value = 8571.143
$ws.Range("I34").Value = 7500
$ws.Range("J34").Value = 9999.333000000001
$ws.Range("K34").Value = 7500
$ws.Range("L34").Value = 9999.333000000001
$ws.Range("M34").Value = -7328
$ws.Range("N34").Value = -10343.333
$ws.Range("H82").Value = 1837.2
$ws.Range("I82").Value = 1760
$ws.Range("K82").Value = 1760
$ws.Range("M82").Value = -1399
$ws.Range("H85").Value = 1837.2
$ws.Range("I85").Value = 1760
$ws.Range("K85").Value = 1760
$ws.Range("M85").Value = -512
$ws.Range("H100").Value = 18528.143
$ws.Range("I100").Value = 2938.6
$ws.Range("K100").Value = 2938.6
$ws.Range("M100").Value = -2397.6

$ws = $wb.Worksheets("WVR")
$ws.Range("H43").Value = 25000
$ws.Range("J43").Value = 25000
$ws.Range("L43").Value = 25000
$ws.Range("N43").Value = -25298
$ws.Range("H46").Value = 82268
$ws.Range("J46").Value = 82268
$ws.Range("L46").Value = 82268
$ws.Range("N46").Value = -82730
$ws.Range("H107").Value = 3697.5833
$ws.Range("I107").Value = 2486
$ws.Range("J107").Value = 5393.8
$ws.Range("K107").Value = 7458
$ws.Range("L107").Value = 16181.4
$ws.Range("M107").Value = -5538
$ws.Range("N107").Value = -20021.4
$ws.Range("H134").Value = 82268
$ws.Range("J134").Value = 82268
$ws.Range("L134").Value = 246804
$ws.Range("N134").Value = -251874
